$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.359.51"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.566.43"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3735"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.13"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3359"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07421"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.116"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.852"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.839"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").Value = "1.565.88"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001106"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06686"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.118"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("D24").Value = "22.350.19"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.368"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.501"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -10.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").Value = "1.738.92"
$ws.Range("E31").Value = "  +0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.986"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9745"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.876"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.635"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08400"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.381"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02449"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2249"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06343"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.333"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6166"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.78"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5739"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.029"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.19"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07288"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.21%  "

$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +0.49%  "
